# "Createdata" refresh: the old static 'UserDetails' demo sheet (Name/Country
# sample rows) is retired, and the 'Sheet1' sheet - which already holds the
# UserDetails column headers (firstName/lastName/userName/password/email/
# mobileNumber) - is renamed to 'UserDetails' and populated with one freshly
# generated user record ("UpdatedDataInExcel").

$wb = $excel.ActiveWorkbook

# Drop the old sample-data sheet.
$oldUserDetails = $wb.Worksheets.Item("UserDetails")
$oldUserDetails.Delete()

# The remaining sheet (header-only) becomes the new UserDetails sheet.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "UserDetails"
$ws.Activate()

# Write the generated record - field order matches the Createdata method's
# generation order (firstName, lastName, password, mobileNumber, then the
# derived userName and the email), not strictly left-to-right by column.
$ws.Cells.Item(2, 1).Value = "Crystal"
$ws.Cells.Item(2, 2).Value = "Zulauf"
$ws.Cells.Item(2, 4).Value = "up0t5s8v1"
$ws.Cells.Item(2, 6).Value = "1-626-558-2900"
$ws.Cells.Item(2, 3).Value = "CrystalZulauf35660"
$ws.Cells.Item(2, 5).Value = "ignacio.bernhard@yahoo.com"

$ws.Range("F2").Select()
